$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column for rows 2-5 from 2023-09-14 (45183)
# to 2023-09-15 (45184), keeping the existing date serial numbering.
foreach ($row in 2..5) {
    $ws.Cells.Item($row, 3).Value = 45184
}
